$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Copy formatting from the last existing data row (118) down into the two
# new rows (119, 120) so the new cells pick up the same styles (e.g. the
# date number format used by columns C and D).
$ws.Range("A118:I118").Copy()
$ws.Range("A119:I120").PasteSpecial(-4122)

# Row 119 - American Moonshot
$ws.Cells.Item(119, 1).Value = "American Moonshot"
$ws.Cells.Item(119, 2).Value = "Douglas Brinkley"
$ws.Cells.Item(119, 3).Value = 44053
$ws.Cells.Item(119, 4).Value = 44056
$ws.Cells.Item(119, 5).Value = "space race;cold war;history;john f kennedy;moon landing"
$ws.Cells.Item(119, 6).Value = "Audio"
$ws.Cells.Item(119, 7).Value = "17 Hours 23 Mins"
$ws.Cells.Item(119, 8).Value = 3
$ws.Cells.Item(119, 9).Value = $true

# Row 120 - Lords of Finance
$ws.Cells.Item(120, 1).Value = "Lords of Finance"
$ws.Cells.Item(120, 2).Value = "Liaquat Ahamed"
$ws.Cells.Item(120, 3).Value = 44034
$ws.Cells.Item(120, 4).Value = 44057
$ws.Cells.Item(120, 5).Value = "history;world war 1;world war 2;finance;reparations;great depression"
$ws.Cells.Item(120, 6).Value = "Hard Copy"
$ws.Cells.Item(120, 7).Value = "505 Pages"
$ws.Cells.Item(120, 8).Value = 3
$ws.Cells.Item(120, 9).Value = $true

# Update the view so the new last row is visible/selected, matching the
# author's workbook state after adding the rows.
$ws.Range("A121").Select()
